{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the last paragraph (\"Predicting Fingers\") to insert the new\n// content directly after it.\nconst last = paragraphs.getLast();\nlast.load(\"text\");\nawait context.sync();\n\nconst empty = last.insertParagraph(\"\", \"After\");\nconst heading = empty.insertParagraph(\"1) Define the Problem\", \"After\");\nheading.insertParagraph(\n  \"The key to this problem is figuring out the mathematical pattern so that you do not have to literally count to 1000 in order to find out what finger the girl would stop on when counting to 1000, or any other number.\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# \"Predicting Fingers\" is currently the last paragraph in the document\n# body. Append three new paragraphs after it: a blank spacer, the new\n# \"1) Define the Problem\" heading, and the explanatory text for it.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n\n$count = $d.Paragraphs.Count\n\n$headingPara = $d.Paragraphs.Item($count - 1)\n$headingPara.Range.Text = \"1) Define the Problem\"\n\n$bodyPara = $d.Paragraphs.Item($count)\n$bodyPara.Range.Text = \"The key to this problem is figuring out the mathematical pattern so that you do not have to literally count to 1000 in order to find out what finger the girl would stop on when counting to 1000, or any other number.\"\n"}
